$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows that are no longer present in the updated dataset (self-pairs / dropped combos).
# Delete from the bottom up so row indices of earlier rows remain valid.
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(2).Delete()

# Write the refreshed TPM-derived values into the remaining rows (now rows 2-9).
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna1"
$ws.Range("C2").Value = "Epha5"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 23.73148533333334
$ws.Range("H2").Value = 71.194456
$ws.Range("I2").Value = 0.8653076146801144
$ws.Range("J2").Value = 0.8653076146801145
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06103333333333334
$ws.Range("N2").Value = 0.1831
$ws.Range("O2").Value = 0.2094151016766933
$ws.Range("P2").Value = 0.2094151016766933
$ws.Range("Q2").Value = 1.448411654844445
$ws.Range("R2").Value = 13.0357048936
$ws.Range("S2").Value = 0.1812084821098531
$ws.Range("T2").Value = 0.1812084821098531

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna1"
$ws.Range("C3").Value = "Epha5"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 23.73148533333334
$ws.Range("H3").Value = 71.194456
$ws.Range("I3").Value = 0.8653076146801144
$ws.Range("J3").Value = 0.8653076146801145
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2304133333333333
$ws.Range("N3").Value = 0.69124
$ws.Range("O3").Value = 0.7905848983233067
$ws.Range("P3").Value = 0.7905848983233067
$ws.Range("Q3").Value = 5.468050640604445
$ws.Range("R3").Value = 49.21245576544
$ws.Range("S3").Value = 0.6840991325702613
$ws.Range("T3").Value = 0.6840991325702613

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efna1"
$ws.Range("C4").Value = "Epha5"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.471212
$ws.Range("H4").Value = 7.413636
$ws.Range("I4").Value = 0.09010639372350319
$ws.Range("J4").Value = 0.09010639372350321
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.06103333333333334
$ws.Range("N4").Value = 0.1831
$ws.Range("O4").Value = 0.2094151016766933
$ws.Range("P4").Value = 0.2094151016766933
$ws.Range("Q4").Value = 0.1508263057333333
$ws.Range("R4").Value = 1.3574367516
$ws.Range("S4").Value = 0.01886963960332758
$ws.Range("T4").Value = 0.01886963960332758

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha5"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.471212
$ws.Range("H5").Value = 7.413636
$ws.Range("I5").Value = 0.09010639372350319
$ws.Range("J5").Value = 0.09010639372350321
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2304133333333333
$ws.Range("N5").Value = 0.69124
$ws.Range("O5").Value = 0.7905848983233067
$ws.Range("P5").Value = 0.7905848983233067
$ws.Range("Q5").Value = 0.5694001942933333
$ws.Range("R5").Value = 5.12460174864
$ws.Range("S5").Value = 0.07123675412017562
$ws.Range("T5").Value = 0.07123675412017562

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Efna1"
$ws.Range("C6").Value = "Epha5"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.148663
$ws.Range("H6").Value = 3.445989
$ws.Range("I6").Value = 0.04188304383987305
$ws.Range("J6").Value = 0.04188304383987305
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.004310666666666667
$ws.Range("N6").Value = 0.012932
$ws.Range("O6").Value = 0.2094151016766933
$ws.Range("P6").Value = 0.2094151016766933
$ws.Range("Q6").Value = 0.07010673176666667
$ws.Range("R6").Value = 0.6309605859
$ws.Range("S6").Value = 0.008770941884256417
$ws.Range("T6").Value = 0.008770941884256417

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Efna1"
$ws.Range("C7").Value = "Epha5"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.148663
$ws.Range("H7").Value = 3.445989
$ws.Range("I7").Value = 0.04188304383987305
$ws.Range("J7").Value = 0.04188304383987305
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2304133333333333
$ws.Range("N7").Value = 0.69124
$ws.Range("O7").Value = 0.7905848983233067
$ws.Range("P7").Value = 0.7905848983233067
$ws.Range("Q7").Value = 0.2646672707066667
$ws.Range("R7").Value = 2.38200543636
$ws.Range("S7").Value = 0.03311210195561663
$ws.Range("T7").Value = 0.03311210195561663

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Efna1"
$ws.Range("C8").Value = "Epha5"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.07412966666666666
$ws.Range("H8").Value = 0.222389
$ws.Range("I8").Value = 0.002702947756509242
$ws.Range("J8").Value = 0.002702947756509243
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.06103333333333334
$ws.Range("N8").Value = 0.1831
$ws.Range("O8").Value = 0.2094151016766933
$ws.Range("P8").Value = 0.2094151016766933
$ws.Range("Q8").Value = 0.004524380655555556
$ws.Range("R8").Value = 0.0407194259
$ws.Range("S8").Value = 0.000566038079256173
$ws.Range("T8").Value = 0.0005660380792561731

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Efna1"
$ws.Range("C9").Value = "Epha5"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.07412966666666666
$ws.Range("H9").Value = 0.222389
$ws.Range("I9").Value = 0.002702947756509242
$ws.Range("J9").Value = 0.002702947756509243
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2304133333333333
$ws.Range("N9").Value = 0.69124
$ws.Range("O9").Value = 0.7905848983233067
$ws.Range("P9").Value = 0.7905848983233067
$ws.Range("Q9").Value = 0.01708046359555555
$ws.Range("R9").Value = 0.15372417236
$ws.Range("S9").Value = 0.002136909677253069
$ws.Range("T9").Value = 0.002136909677253069
